$d = $word.ActiveDocument

# --- Step 1: "harbon" + "</m>" + " " -> "harbon " -----------------------
# "harbon" is a unique run in the doc. We can't just Find/Replace across the
# three runs directly (that would coalesce the preceding "c" run into the
# replacement too), so instead: locate "harbon", collapse to its end, grow
# the (now tiny) range just far enough to cover the following "</m>" run,
# and delete *that* text in a tightly scoped Find so only the run adjacent
# to "harbon" is touched. "harbon" then naturally absorbs the following
# (already-matching-format) space run into "harbon ".
$rng1 = $d.Content
$rng1.Find.Execute("harbon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(0)          # wdCollapseEnd
$rng1.MoveEnd(1, 4)        # wdCharacter - grow just enough to cover "</m>"
$rng1.Find.Execute("</m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Step 2: "</tl>" -> "</m></tl>" (only the one after "allume") -------
$rng2 = $d.Content
$rng2.Find.Execute("allume", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)          # wdCollapseEnd
$rng2.MoveEnd(1, 5)        # grow just enough to cover "</tl>"
$rng2.Find.Execute("</tl>", $true, $false, $false, $false, $false, $true, 1, $false, "</m></tl>", 2) | Out-Null
